# CRC Reporting Rates — add the 2024 Winter totals (row 19, cols D & E)
# and move the sheet's active-cell selection from D21 to E21, matching
# the author's re-upload of the parameter workbook with the 2024
# harvest-estimate inputs filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 19 is year 2024: A19=2024, B19=187214, C19=81091 already present.
# Fill in the newly-reported Winter totals.
$ws.Range("D19").Value = 33743
$ws.Range("E19").Value = 16660

# The workbook was left with E21 selected (was D21) when last saved.
[void]$ws.Range("E21").Select()
